# Reorder the "Recorded By" (column G) values so that
# "dnasr281@gmail.com" is listed first, ahead of "System" / "admin@admin.com".
#
# Example: "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
#          "admin@admin.com, dnasr281@gmail.com" -> "dnasr281@gmail.com, admin@admin.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Value()

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "admin@admin.com, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, admin@admin.com"
    }
}
